# Auto-generated script applying market-data refresh values
# produced by the scheduled runner, matching the target OOXML diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 301.6
$ws.Range("I4").Value = 300.91666
$ws.Range("K4").Value = 300.91666
$ws.Range("M4").Value = -186.91666
$ws.Range("H5").Value = 181.15384
$ws.Range("I5").Value = 70.75
$ws.Range("J5").Value = 357.8
$ws.Range("K5").Value = 70.75
$ws.Range("L5").Value = 357.8
$ws.Range("M5").Value = 44.25
$ws.Range("N5").Value = -587.8
$ws.Range("H17").Value = 1346.6666
$ws.Range("J17").Value = 1570
$ws.Range("L17").Value = 4710
$ws.Range("N17").Value = -5046
$ws.Range("H40").Value = 1450
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 1450
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 1450
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -1800
$ws.Range("H86").Value = 6524.5
$ws.Range("I86").Value = 4649
$ws.Range("K86").Value = 4649
$ws.Range("M86").Value = -3526
$ws.Range("H89").Value = 6524.5
$ws.Range("I89").Value = 4649
$ws.Range("K89").Value = 23245
$ws.Range("M89").Value = -17629
$ws.Range("H100").Value = 2600.3333
$ws.Range("I100").Value = 2600.3333
$ws.Range("K100").Value = 2600.3333
$ws.Range("M100").Value = -2059.3333
$ws.Range("H101").Value = 3725
$ws.Range("I101").Value = 6200
$ws.Range("J101").Value = 1250
$ws.Range("K101").Value = 18600
$ws.Range("L101").Value = 3750
$ws.Range("M101").Value = -16978
$ws.Range("N101").Value = -6994
$ws.Range("H112").Value = 3125.9443
$ws.Range("J112").Value = 3227.5881
$ws.Range("L112").Value = 9682.764299999999
$ws.Range("N112").Value = -11898.7643
$ws.Range("H137").Value = 1968
$ws.Range("I137").Value = 1829.5
$ws.Range("K137").Value = 5488.5
$ws.Range("M137").Value = -2938.5
$ws.Range("H138").Value = 2663.25
$ws.Range("I138").Value = 1521.9375
$ws.Range("J138").Value = 4185
$ws.Range("K138").Value = 4565.8125
$ws.Range("L138").Value = 12555
$ws.Range("M138").Value = 574.1875
$ws.Range("N138").Value = -22835

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 14766
$ws.Range("I28").Value = 14766
$ws.Range("K28").Value = 14766
$ws.Range("M28").Value = -14574
$ws.Range("H32").Value = 4710.647
$ws.Range("I32").Value = 3607.9285
$ws.Range("K32").Value = 3607.9285
$ws.Range("M32").Value = -3320.9285
$ws.Range("H63").Value = 6133.3335
$ws.Range("I63").Value = 1900
$ws.Range("J63").Value = 8250
$ws.Range("K63").Value = 1900
$ws.Range("L63").Value = 8250
$ws.Range("M63").Value = -1214
$ws.Range("N63").Value = -9622
$ws.Range("H66").Value = 6133.3335
$ws.Range("I66").Value = 1900
$ws.Range("J66").Value = 8250
$ws.Range("K66").Value = 9500
$ws.Range("L66").Value = 41250
$ws.Range("M66").Value = -6068
$ws.Range("N66").Value = -48114
$ws.Range("H97").Value = 1743.1
$ws.Range("J97").Value = 2201.4
$ws.Range("L97").Value = 2201.4
$ws.Range("N97").Value = -3193.4
$ws.Range("H99").Value = 14766
$ws.Range("I99").Value = 14766
$ws.Range("K99").Value = 14766
$ws.Range("M99").Value = -11771

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3758.5
$ws.Range("I20").Value = 2385
$ws.Range("K20").Value = 2385
$ws.Range("M20").Value = -2138
$ws.Range("H86").Value = 2912.0908
$ws.Range("I86").Value = 670.3333
$ws.Range("J86").Value = 13000
$ws.Range("K86").Value = 670.3333
$ws.Range("L86").Value = 13000
$ws.Range("M86").Value = 452.6667
$ws.Range("N86").Value = -15246
$ws.Range("H89").Value = 2912.0908
$ws.Range("I89").Value = 670.3333
$ws.Range("J89").Value = 13000
$ws.Range("K89").Value = 3351.6665
$ws.Range("L89").Value = 65000
$ws.Range("M89").Value = 2264.3335
$ws.Range("N89").Value = -76232
$ws.Range("H94").Value = 2475.5
$ws.Range("I94").Value = 1213.3334
$ws.Range("J94").Value = 4999.8335
$ws.Range("K94").Value = 1213.3334
$ws.Range("L94").Value = 4999.8335
$ws.Range("M94").Value = -762.3334
$ws.Range("N94").Value = -5901.8335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1600.5714
$ws.Range("I31").Value = 1582.5454
$ws.Range("J31").Value = 1666.6666
$ws.Range("K31").Value = 1582.5454
$ws.Range("L31").Value = 1666.6666
$ws.Range("M31").Value = -1287.5454
$ws.Range("N31").Value = -2256.6666
$ws.Range("H33").Value = 2976
$ws.Range("I33").Value = 2976
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 2976
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -2597
$ws.Range("N33").ClearContents()
$ws.Range("H34").Value = 1600.5714
$ws.Range("I34").Value = 1582.5454
$ws.Range("J34").Value = 1666.6666
$ws.Range("K34").Value = 1582.5454
$ws.Range("L34").Value = 1666.6666
$ws.Range("M34").Value = -1380.5454
$ws.Range("N34").Value = -2070.6666
$ws.Range("H107").Value = 688
$ws.Range("H132").Value = 2956.2068
$ws.Range("I132").Value = 2943.7917
$ws.Range("K132").Value = 8831.375100000001
$ws.Range("M132").Value = -6301.375100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 680.6
$ws.Range("I86").Value = 966.6667
$ws.Range("J86").Value = 251.5
$ws.Range("K86").Value = 2900.0001
$ws.Range("L86").Value = 754.5
$ws.Range("M86").Value = -1714.0001
$ws.Range("N86").Value = -3126.5
$ws.Range("H89").Value = 680.6
$ws.Range("I89").Value = 966.6667
$ws.Range("J89").Value = 251.5
$ws.Range("K89").Value = 8700.0003
$ws.Range("L89").Value = 2263.5
$ws.Range("M89").Value = -2772.0003
$ws.Range("N89").Value = -14119.5
$ws.Range("H122").Value = 451.625
$ws.Range("I122").Value = 117.2
$ws.Range("K122").Value = 1054.8
$ws.Range("M122").Value = 1395.2
$ws.Range("H132").Value = 2069.2856
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 2497
$ws.Range("K132").Value = 9000
$ws.Range("L132").Value = 22473
$ws.Range("M132").Value = -6470
$ws.Range("N132").Value = -27533

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 23691
$ws.Range("J98").Value = 23691
$ws.Range("L98").Value = 23691
$ws.Range("N98").Value = -29681
$ws.Range("H107").Value = 9813.5
$ws.Range("I107").Value = 1995
$ws.Range("J107").Value = 13722.75
$ws.Range("K107").Value = 1995
$ws.Range("L107").Value = 13722.75
$ws.Range("M107").Value = -75
$ws.Range("N107").Value = -17562.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3002.5
$ws.Range("I100").Value = 3002.5
$ws.Range("K100").Value = 3002.5
$ws.Range("M100").Value = -2461.5
$ws.Range("H122").Value = 5000
$ws.Range("I122").Value = 5000
$ws.Range("K122").Value = 15000
$ws.Range("M122").Value = -12550
$ws.Range("H132").Value = 1637.3077
$ws.Range("I132").Value = 1565.4166
$ws.Range("K132").Value = 4696.2498
$ws.Range("M132").Value = -2166.2498

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2319.6
$ws.Range("I122").Value = 2319.6
$ws.Range("K122").Value = 6958.799999999999
$ws.Range("M122").Value = -4508.799999999999
$ws.Range("H132").Value = 1704.7222
$ws.Range("I132").Value = 1755
$ws.Range("K132").Value = 5265
$ws.Range("M132").Value = -2735

